# Update "想去人数" (column F) counts for several events that appear in both
# the "展览" sheet and the "全部类型" sheet (which duplicates the same rows).
$wb = $excel.ActiveWorkbook

# Row number -> new value for column F
$updates = @{
    3  = 231
    6  = 12
    7  = 5730
    8  = 5074
    9  = 24
    13 = 212
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
